# Apply crypto price/volume updates from the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.030.11'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '1.955.46'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06830'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('E10').Value = '  -0.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.11%  '
$ws.Range('D12').Value = '1.957.24'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07824'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.460'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7016'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '283.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.36%  '
$ws.Range('D17').Value = '31.061.27'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007700'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '2.202.41'
$ws.Range('E20').Value = '  -0.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.493'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.485'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.830'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.197'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1057'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.606'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.450'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04936'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7646'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.172'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.728'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02010'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.517'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.110'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.96'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.90%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4468'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8856'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.205'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.12%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '1.002.21'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1260'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.365'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2588'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.33%  '
